# Update stats for 2025-07 (row 20 in sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6149
$ws.Range("C20").Value = 976
$ws.Range("D20").Value = 5556174
$ws.Range("E20").Value = 903.5898520084567
$ws.Range("F20").Value = 6.218690620141643
$ws.Range("G20").Value = 3.609341825902335
$ws.Range("H20").Value = 25.68615683630071
